# Enabled Dock frames to clear out when no content available from selected doc
#
# The "frames" table in column H (docs/ID/ptr/name) gets a new blank line
# above it so it lines up with a new "userType" row inserted into the
# adjoining F column table - mirroring the existing userType/docID pair
# higher up in the sheet. Net effect: H7:H10 shift down to H8:H11, F8
# ("docID") shifts down to F9 and is replaced at F8 by "userType". The
# connector shape that pointed at the old F7/H7 row moves down with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the "docs / ID / ptr / name" block in column H down by one row
# (H10->H11, H9->H10, H8->H9, H7->H8), freeing up H7.
$ws.Range("H11").Value = $ws.Range("H10").Value()
$ws.Range("H10").Value = $ws.Range("H9").Value()
$ws.Range("H9").Value = $ws.Range("H8").Value()
$ws.Range("H8").Value = $ws.Range("H7").Value()
$ws.Range("H7").ClearContents()

# Push "docID" down to F9 and put the new "userType" label at F8.
$ws.Range("F9").Value = $ws.Range("F8").Value()
$ws.Range("F8").Value = "userType"

# The connector shape anchored at row 7/8 (0-indexed) needs to follow the
# row-8 content down to row 8/9 (0-indexed), i.e. one row height (14.4pt).
$shp = $ws.Shapes.Item("Straight Arrow Connector 21")
$shp.Top = 122.00456692913386

# New active selection left behind by the edit.
$ws.Range("F11").Select()
